# Adds the hypothesis statements (H0/H1) and the reject/cannot-reject
# conclusion to every "Rocket vs X" sheet, matching the author's commit:
# "added hypothesis's and the conclusions of the tests on whether or not
#  I can reject the null hypothesis"
#
# Layout on every sheet:
#   row 17  (merged A17:H17)  -> H0 statement
#   row 18  (merged A18:H18)  -> H1 statement
#   row 20  (merged A20:F20)  -> conclusion (reject / cannot reject)
# Sheet 5 ("Rocket vs kNN-DTW") additionally already has data in columns
# H:J for rows 16-18, so row 19 is left as a blank spacer row there.

$wb = $excel.ActiveWorkbook

$xlCenter = -4108
$xlLineStyleNone = -4142

$CANNOT_REJECT = "Based off the evidence of this test we cannot reject the null hypothesis "
$CAN_REJECT = "Based off the evidence of this test we can reject the null hypothesis "

function Add-HypothesisRows {
    param(
        $ws,
        [string]$h0,
        [string]$h1,
        [string]$conclusion
    )

    $ws.Range("A17:H17").HorizontalAlignment = $xlCenter
    $ws.Range("A17").Value = $h0
    $ws.Range("A17:H17").Merge()

    $ws.Range("A18:H18").HorizontalAlignment = $xlCenter
    $ws.Range("A18").Value = $h1
    $ws.Range("A18:H18").Merge()

    $ws.Range("A20:F20").HorizontalAlignment = $xlCenter
    $ws.Range("A20").Value = $conclusion
    $ws.Range("A20:F20").Merge()
}

# --- Sheet 1: Rocket vs AdaBoost -------------------------------------------
$s1 = $wb.Worksheets.Item(1)

# --- Sheet 2: Rocket vs Boss Ensemble ---------------------------------------
$s2 = $wb.Worksheets.Item(2)

# Hypotheses for sheet 1 and sheet 2 are authored before either conclusion,
# then sheet 2's "cannot reject" conclusion is authored first, followed by
# sheet 1's "can reject" conclusion -- this ordering reproduces the exact
# shared-string table order of the target workbook.
$s1.Range("A17:H17").HorizontalAlignment = $xlCenter
$s1.Range("A17").Value = "H0: there is no difference in mean accuracry between Rocket and AdaBoost on gym movements"
$s1.Range("A17:H17").Merge()

$s1.Range("A18:H18").HorizontalAlignment = $xlCenter
$s1.Range("A18").Value = "H1: there is a difference in mean accuracry between Rocket and AdaBoost on gym movements"
$s1.Range("A18:H18").Merge()

$s2.Range("A17:H17").HorizontalAlignment = $xlCenter
$s2.Range("A17").Value = "H0: there is no difference in mean accuracry between Rocket and Boss Ensemble on gym movements"
$s2.Range("A17:H17").Merge()

$s2.Range("A18:H18").HorizontalAlignment = $xlCenter
$s2.Range("A18").Value = "H1: there is a difference in mean accuracry between Rocket and Boss Ensemble on gym movements"
$s2.Range("A18:H18").Merge()

$s2.Range("A20:F20").HorizontalAlignment = $xlCenter
$s2.Range("A20").Value = $CANNOT_REJECT
$s2.Range("A20:F20").Merge()

$s1.Range("A20:F20").HorizontalAlignment = $xlCenter
$s1.Range("A20").Value = $CAN_REJECT
$s1.Range("A20:F20").Merge()

# --- Sheet 3: Rocket vs Decision Tree ---------------------------------------
$s3 = $wb.Worksheets.Item(3)
Add-HypothesisRows $s3 `
    "H0: there is no difference in mean accuracry between Rocket and Decision Tree on gym movements" `
    "H1: there is a difference in mean accuracry between Rocket and Decision Tree on gym movements" `
    $CAN_REJECT

# --- Sheet 4: Rocket vs kNN-ED ----------------------------------------------
$s4 = $wb.Worksheets.Item(4)
Add-HypothesisRows $s4 `
    "H0: there is no difference in mean accuracry between Rocket and kNN-ED on gym movements" `
    "H1: there is a difference in mean accuracry between Rocket and kNN-ED on gym movements" `
    $CAN_REJECT

# --- Sheet 5: Rocket vs kNN-DTW ----------------------------------------------
# This sheet already has content in H:J for rows 16-18, so the new rows keep
# that content and row 19 is an untouched (but format-probed) blank row.
$s5 = $wb.Worksheets.Item(5)

$s5.Range("A17:H17").HorizontalAlignment = $xlCenter
$s5.Range("A17").Value = "H0: there is no difference in mean accuracry between Rocket and kNN-DTW on gym movements"
$s5.Range("A17:H17").Merge()

$s5.Range("A18:H18").HorizontalAlignment = $xlCenter
$s5.Range("A18").Value = "H1: there is a difference in mean accuracry between Rocket and kNN-DTW on gym movements"
$s5.Range("A18:H18").Merge()

$s5.Range("A19:G19").Borders.LineStyle = $xlLineStyleNone

$s5.Range("A20:F20").HorizontalAlignment = $xlCenter
$s5.Range("A20").Value = $CANNOT_REJECT
$s5.Range("A20:F20").Merge()
$s5.Range("G20").Borders.LineStyle = $xlLineStyleNone

# --- Sheet 6: Rocket vs MLP --------------------------------------------------
$s6 = $wb.Worksheets.Item(6)
Add-HypothesisRows $s6 `
    "H0: there is no difference in mean accuracry between Rocket and MLP on gym movements" `
    "H1: there is a difference in mean accuracry between Rocket and MLP on gym movements" `
    $CAN_REJECT

# --- Sheet 7: Rocket vs Naive Bayes ------------------------------------------
$s7 = $wb.Worksheets.Item(7)
Add-HypothesisRows $s7 `
    "H0: there is no difference in mean accuracry between Rocket and Naïve Bayes on gym movements" `
    "H1: there is a difference in mean accuracry between Rocket and Naïve Bayes on gym movements" `
    $CAN_REJECT

# --- Sheet 8: Rocket vs Random Forest ----------------------------------------
$s8 = $wb.Worksheets.Item(8)
Add-HypothesisRows $s8 `
    "H0: there is no difference in mean accuracry between Rocket and Random Forest on gym movements" `
    "H1: there is a difference in mean accuracry between Rocket and Random Forest on gym movements" `
    $CAN_REJECT

# --- Sheet 9: Rocket vs Time Series Forest -----------------------------------
$s9 = $wb.Worksheets.Item(9)
Add-HypothesisRows $s9 `
    "H0: there is no difference in mean accuracry between Rocket and Time Series Forest on gym movements" `
    "H1: there is a difference in mean accuracry between Rocket and Time Series Forest on gym movements" `
    $CAN_REJECT

# --- Restore each sheet's on-screen selection -------------------------------
# (Selecting a range also makes that sheet the active tab, so sheet 9 -- the
# sheet that was active before the edit -- is selected last to keep it active.)
$s1.Range("A17:H20").Select() | Out-Null
$s2.Range("A17:H20").Select() | Out-Null
$s3.Range("A17:H20").Select() | Out-Null
$s4.Range("A17:H20").Select() | Out-Null
$s5.Range("A17:H20").Select() | Out-Null
$s6.Range("A17:H20").Select() | Out-Null
$s7.Range("A17:H20").Select() | Out-Null
$s8.Range("A17:H20").Select() | Out-Null
$s9.Range("H23").Select() | Out-Null
